$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay as text, matching the
# source data which stores formatted/locale strings (e.g. "1.004", "25.908.94")
# rather than numeric values.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.908.94'
$ws.Range("E2").Value = '  -0.30%  '

$ws.Range("D3").Value = '1.638.60'
$ws.Range("E3").Value = '  -0.36%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.29%  '

$ws.Range("D5").Value = '214.67'
$ws.Range("E5").Value = '  -0.34%  '

$ws.Range("D6").Value = '0.5062'
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").Value = '0.2552'

$ws.Range("E9").Value = '  -0.82%  '

$ws.Range("D10").Value = '19.46'
$ws.Range("E10").Value = '  -1.22%  '

$ws.Range("D11").Value = '0.07748'
$ws.Range("E11").Value = '  -0.04%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.648.37'
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.271'
$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("D14").Value = '0.5433'
$ws.Range("E14").Value = '  -0.66%  '

$ws.Range("D15").Value = '0.0₅7810'
$ws.Range("E15").Value = '  -1.75%  '

$ws.Range("D16").Value = '64.13'
$ws.Range("E16").Value = '  -0.49%  '

$ws.Range("D17").Value = '25.944.56'
$ws.Range("E17").Value = '  -0.19%  '

$ws.Range("D18").Value = '1.004'
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").Value = '196.25'
$ws.Range("E19").Value = '  -2.90%  '

$ws.Range("D20").Value = '4.448'
$ws.Range("E20").Value = '  +1.29%  '

$ws.Range("D21").Value = '9.922'
$ws.Range("E21").Value = '  +0.18%  '

$ws.Range("D22").Value = '6.011'
$ws.Range("E22").Value = '  +0.29%  '

$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  +0.23%  '

$ws.Range("D24").Value = '1.888'
$ws.Range("E24").Value = '  +0.70%  '

$ws.Range("D25").Value = '140.98'
$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("D26").Value = '0.1185'
$ws.Range("E26").Value = '  +3.86%  '

$ws.Range("D27").Value = '6.862'
$ws.Range("E27").Value = '  +0.51%  '

$ws.Range("E28").Value = '  -0.06%  '

$ws.Range("E29").Value = '  -0.51%  '

$ws.Range("D30").Value = '0.04932'
$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("D31").Value = '3.250'
$ws.Range("E31").Value = '  -0.70%  '

$ws.Range("D32").Value = '3.172'
$ws.Range("E32").Value = '  -1.30%  '

$ws.Range("D33").Value = '1.536'
$ws.Range("E33").Value = '  -0.64%  '

$ws.Range("D34").Value = '2.367'
$ws.Range("E34").Value = '  +0.17%  '

$ws.Range("D35").Value = '0.8926'
$ws.Range("E35").Value = '  -0.24%  '

$ws.Range("D36").Value = '2.582'
$ws.Range("E36").Value = '  -1.44%  '

$ws.Range("D37").Value = '1.130.91'
$ws.Range("E37").Value = '  -1.90%  '

$ws.Range("D38").Value = '0.5425'
$ws.Range("E38").Value = '  -2.98%  '

$ws.Range("D39").Value = '0.01553'
$ws.Range("E39").Value = '  -0.93%  '

$ws.Range("E40").Value = '  +0.21%  '

$ws.Range("D41").Value = '2.545'
$ws.Range("E41").Value = '  -0.58%  '

$ws.Range("D42").Value = '0.0₈128'
$ws.Range("E42").Value = '  +8.68%  '

$ws.Range("D43").Value = '5.573'
$ws.Range("E43").Value = '  -2.36%  '

$ws.Range("D44").Value = '0.8137'
$ws.Range("E44").Value = '  +0.59%  '

$ws.Range("D45").Value = '99.32'
$ws.Range("E45").Value = '  -0.44%  '

$ws.Range("D46").Value = '1.775.97'
$ws.Range("E46").Value = '  -0.30%  '

$ws.Range("D47").Value = '0.4541'
$ws.Range("E47").Value = '  +0.51%  '

$ws.Range("D48").Value = '1.002'
$ws.Range("E48").Value = '  -0.06%  '

$ws.Range("D49").Value = '54.71'
$ws.Range("E49").Value = '  -0.17%  '

$ws.Range("D50").Value = '0.05072'
$ws.Range("E50").Value = '  +0.54%  '

$ws.Range("D51").Value = '1.006'
$ws.Range("E51").Value = '  +0.28%  '
